$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.976.68'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.13%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.971.61'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.45%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '526.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.31'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.67%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.967.34'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.43%  '
$ws.Range("E9").Value = '  -1.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.24'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.24%  '
$ws.Range("E11").Value = '  -1.59%  '
$ws.Range("E12").Value = '  -2.11%  '
$ws.Range("E13").Value = '  -1.97%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.16'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.457.91'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.99%  '
$ws.Range("E16").Value = '  +0.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.973.81'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.966.87'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.43%  '
$ws.Range("E19").Value = '  -0.68%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '455.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.82%  '
$ws.Range("E21").Value = '  -0.22%  '
$ws.Range("E22").Value = '  -2.98%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.78'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.75%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '77.83'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.68'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.94%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("E27").Value = '  -1.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.59'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("E30").Value = '  -0.88%  '
$ws.Range("B31").Value = 'Mantle'
$ws.Range("C31").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.12'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.48%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.82'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.43%  '
$ws.Range("B33").Value = 'Stacks'
$ws.Range("C33").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.24'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.01%  '
$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '54.20'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.30'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.88%  '
$ws.Range("E36").Value = '  -1.53%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '450.06'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.135.44'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.19%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0379'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.62%  '
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0772'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.37%  '
$ws.Range("E41").Value = '  +5.26%  '
$ws.Range("E42").Value = '  -0.14%  '
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.37'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.37%  '
$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.10%  '
$ws.Range("E45").Value = '  -1.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '24.96'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.50%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '120.18'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.35%  '
$ws.Range("E48").Value = '  +0.73%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₃0499'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.38%  '
$ws.Range("E51").Value = '  +6.44%  '
